# Marks-service doc update:
#   - "Used message resources" paragraph gains a second run " for Labels"
#   - "Used CrudRepository" paragraph becomes "Used message resources" + " for Errors"
#     (keeping the _GoBack bookmark in place)
#   - a brand-new paragraph "Used CrudRepository" is added right after it
#
# Starting layout (Paragraphs 3 & 4):
#   3: "Used message resources"
#   4: "Used CrudRepository" [+ bookmarkStart/bookmarkEnd _GoBack]

$d = $word.ActiveDocument

# --- Step 1: make room for the new trailing "Used CrudRepository" paragraph
# right after the bookmarked paragraph (paragraph 4).
$pCrud = $d.Paragraphs.Item(4)
$pCrud.Range.InsertParagraphAfter()

$pNew = $d.Paragraphs.Item(5)
$pNew.Range.Text = "Used CrudRepository"

# --- Step 2: paragraph 4 text changes from "Used CrudRepository" to
# "Used message resources", then gets a new run " for Errors" appended,
# with the bookmark staying on the paragraph.
$pCrud = $d.Paragraphs.Item(4)
[void]$pCrud.Range.Find.Execute("Used CrudRepository", $true, $false, $false, $false, $false, $true, 1, $false, "Used message resources", 2)

$pCrud = $d.Paragraphs.Item(4)
$xmlErrors = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t>Used message resources</w:t></w:r><w:r><w:t xml:space="preserve"> for Errors</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
[void]$pCrud.Range.InsertXML($xmlErrors)

# --- Step 3: paragraph 3 ("Used message resources") gains a second run
# " for Labels".
$pMsg = $d.Paragraphs.Item(3)
$xmlLabels = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t>Used message resources</w:t></w:r><w:r><w:t xml:space="preserve"> for Labels</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
[void]$pMsg.Range.InsertXML($xmlLabels)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output "$i -> [$($d.Paragraphs.Item($i).Range.Text)]"
}
